$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.088150978088379
$ws.Range("B1").Value = 2.155479907989502
$ws.Range("C1").Value = 9.494928359985352
$ws.Range("D1").Value = 1.016748547554016
$ws.Range("E1").Value = 1.145743012428284
